# daily auto push: 2026-02-27 02:51 UTC
# Insert a new data row for 2026/02/27 (金) at row 894, shifting the
# existing rows 894:935 down to 895:936.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 894 (pushes old row 894 -> 895, etc.)
$ws.Rows.Item(894).Insert()

# Populate the newly inserted row with the new record.
# Column A holds a date-like string ("yyyy/mm/dd") that must stay literal
# text (matching the sibling cells), not get auto-converted to a date
# serial number by Excel's smart input parsing. Temporarily force the
# cell to Text format while assigning it, then clear the format again so
# the cell ends up with no special number formatting applied (matching
# the rest of the table).
$ws.Cells.Item(894, 1).NumberFormat = "@"
$ws.Cells.Item(894, 1).Value = "2026/02/27"
$ws.Cells.Item(894, 1).ClearFormats()

$ws.Cells.Item(894, 2).Value = "金"
$ws.Cells.Item(894, 3).Value = 8
$ws.Cells.Item(894, 4).Value = 201
